$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 5.282202243804932
$ws.Range("B1").Value = 5.303390026092529
$ws.Range("C1").Value = 8.454256057739258
$ws.Range("D1").Value = 8.30510139465332
$ws.Range("E1").Value = 3.689008235931396
